$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")
$table = $ws.ListObjects.Item("Table1")
$c1 = $table.ListColumns.Add()
$c2 = $table.ListColumns.Add()
$ws.Range("Q1:R8").Insert(-4161)
Write-Output "inserted"
